$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure values are written as text (avoid Excel auto-converting
# numeric-looking strings like "595.59" or "14.30" into real numbers,
# which would lose formatting such as trailing zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.093.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.514.66"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.59"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.43"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.69%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.44%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +7.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.28"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.436"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.133.55"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.84"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000182"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.144.38"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.525.06"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.33"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.30"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "396.62"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.97"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.38"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.539"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000122"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.24"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.24%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.29"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.46"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.07"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.07"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.41"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.81%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "163.50"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.84%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.91"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.91"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.70"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0746"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.41"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.44"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.62"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.797.12"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.87"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "341.20"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.30%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.66"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.52"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.850"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.88%  "
